{"js": "const paras = context.document.body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\n// Paragraph 1 (index 0): date + title line\nparas.items[0].insertText(\"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 -27.11.24: \u26a1\ufe0f\ud83d\ude80\\u000bThe Illusion of State in State-Space Models\", Word.InsertLocation.replace);\n\nparas.items[1].insertText(\"\u05de\u05d0\u05de\u05e8 \u05d7\u05e9\u05d5\u05d1 \u05d6\u05d4 \u05d1\u05d5\u05d7\u05df \u05d0\u05ea \u05d4\u05de\u05d2\u05d1\u05dc\u05d5\u05ea \u05d4\u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9\u05d5\u05ea \u05e9\u05dc State Space Models \u05d0\u05d5 (SSMs), \u05d0\u05e9\u05e8 \u05e6\u05de\u05d7\u05d5 \u05db\u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d4 \u05d7\u05dc\u05d5\u05e4\u05d9\u05ea \u05dc\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd \u05e2\u05d1\u05d5\u05e8 \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05d2\u05d3\u05d5\u05dc\u05d9\u05dd. \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05d3\u05d2\u05d9\u05de\u05d9\u05dd \u05e9\u05dc\u05de\u05e8\u05d5\u05ea \u05e2\u05d9\u05e6\u05d5\u05d1\u05dd \u05e9\u05e0\u05e8\u05d0\u05d4 Recurrent \u05d5\u05d1\u05e2\u05dc \u05de\u05e6\u05d1 (\u05db\u05dc\u05d5\u05de\u05e8 stateful), \u05dc\u05de\u05e2\u05e9\u05d4 SSMs (\u05db\u05de\u05d5 \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd) \u05de\u05d5\u05d2\u05d1\u05dc\u05d9\u05dd \u05d1\u05d0\u05d5\u05e4\u05df \u05d1\u05e1\u05d9\u05e1\u05d9 \u05d1\u05d9\u05db\u05d5\u05dc\u05ea\u05dd \u05dc\u05d1\u05d8\u05d0 \u05d7\u05d9\u05e9\u05d5\u05d1 \\\"\u05e8\u05e6\u05d9\u05e3\\\", \u05de\u05db\u05d9\u05d5\u05d5\u05df \u05e9\u05d0\u05d9\u05e0\u05dd \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05d7\u05e9\u05d1 \u05d3\u05d1\u05e8 \u05de\u05d7\u05d5\u05e5 \u05dc\u05de\u05d7\u05dc\u05e7\u05ea \u05d4\u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea TC0. \u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05de\u05de\u05d7\u05dc\u05e7\u05ea TC0 \u05de\u05d5\u05d2\u05d3\u05e8\u05d5\u05ea \u05db\u05db\u05d0\u05dc\u05d5 \u05e9\u05e0\u05d9\u05ea\u05df \u05dc\u05d9\u05d9\u05e6\u05d2\u05df \u05e2\u05dd \u05e9\u05e8\u05e9\u05e8\u05d0\u05d5\u05ea \u05d1\u05d5\u05dc\u05d9\u05d0\u05e0\u05d9\u05d5\u05ea \u05d1\u05e1\u05d9\u05e1\u05d9\u05d5\u05ea (\u05d5\u05d7\u05d9\u05e9\u05d5\u05d1\u05d9 \u05e1\u05e3 \u05d5- majority vote) \u05d1\u05e2\u05d5\u05de\u05e7 \u05e1\u05d5\u05e4\u05d9 (\u05dc\u05de\u05e9\u05dc \u05d7\u05d9\u05d1\u05d5\u05e8 \u05e9\u05dc \u05de\u05e1\u05e4\u05e8\u05d9\u05dd, \u05de\u05db\u05e4\u05dc\u05d4 \u05d0\u05d5 \u05de\u05d9\u05d5\u05df \u05e9\u05dc n \u05de\u05e1\u05e4\u05e8\u05d9\u05dd). \u05de\u05d3\u05d5\u05d1\u05e8 \u05d1\u05de\u05d7\u05dc\u05e7\u05d4 \u05d4\u05db\u05d9 \\\"\u05e4\u05e9\u05d5\u05d8\u05d4\\\" \u05d1\u05d4\u05d9\u05e8\u05e8\u05db\u05d9\u05d4 \u05e9\u05dc \u05ea\u05d5\u05e8\u05d4 \u05e1\u05d9\u05d1\u05d5\u05db\u05d9\u05d5\u05ea circuit (\u05db\u05dc\u05d5\u05de\u05e8 circuit complexity).\", Word.InsertLocation.replace);\nparas.items[2].insertText(\"\u05de\u05e9\u05de\u05e2\u05d5\u05ea \u05d4\u05d3\u05d1\u05e8 \u05d4\u05d9\u05d0 \u05e9-SSMs \u05d0\u05d9\u05e0\u05dd \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05e4\u05ea\u05d5\u05e8 \u05d1\u05e2\u05d9\u05d5\u05ea \u05de\u05e1\u05d5\u05d2 permutation composition \u05e9- RNNs \u05d1\u05e2\u05dc\u05d5\u05ea \u05e9\u05db\u05d1\u05d4 \u05d0\u05d7\u05ea \u05de\u05e1\u05d5\u05d2\u05dc\u05d5\u05ea \u05dc\u05e4\u05ea\u05d5\u05e8.\", Word.InsertLocation.replace);\nparas.items[3].insertText(\"\u05ea\u05e8\u05d5\u05de\u05d5\u05ea \u05de\u05e8\u05db\u05d6\u05d9\u05d5\u05ea \u05e9\u05dc \u05d4\u05de\u05d0\u05de\u05e8:\", Word.InsertLocation.replace);\nparas.items[4].insertText(\"\u05e0\u05d9\u05ea\u05d5\u05d7 \u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9:\", Word.InsertLocation.replace);\nparas.items[5].insertText(\"\u05de\u05d5\u05db\u05d9\u05d7 \u05e9\u05d2\u05dd SSMs \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05d9\u05dd \u05d5\u05d2\u05dd SSMs \u05d1\u05e1\u05d2\u05e0\u05d5\u05df Mamba \u05de\u05d5\u05d2\u05d1\u05dc\u05d9\u05dd \u05dc\u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea \u05d7\u05d9\u05e9\u05d5\u05d1\u05d9\u05ea TC0\", Word.InsertLocation.replace);\nparas.items[6].insertText(\"\u05de\u05e8\u05d0\u05d4 \u05e9-SSMs \u05d0\u05d9\u05e0\u05dd \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05e4\u05ea\u05d5\u05e8 \u05d1\u05e2\u05d9\u05d5\u05ea \u05e9\u05dc\u05de\u05d5\u05ea-NC1 (\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05e9\u05e0\u05d9\u05ea\u05df \u05dc\u05d9\u05d9\u05e6\u05d2 \u05d0\u05d5\u05ea\u05df \u05e2\u05dd \u05e4\u05e2\u05d5\u05dc\u05d5\u05ea \u05d1\u05d5\u05dc\u05d9\u05d0\u05e0\u05d9\u05d5\u05ea \u05d1\u05e2\u05d5\u05de\u05e7 \u05dc\u05d5\u05d2\u05e8\u05d9\u05ea\u05de\u05d9 \u05de\u05de\u05d9\u05de\u05d3 \u05d4\u05d1\u05e2\u05d9\u05d4 - \u05de\u05e1\u05e4\u05e8 \u05de\u05e9\u05ea\u05e0\u05d9\u05dd \u05d1\u05d2\u05d3\u05d5\u05dc) \u05db\u05de\u05d5 \u05d4\u05e8\u05db\u05d1\u05ea \u05ea\u05de\u05d5\u05e8\u05d5\u05ea. \u05db\u05dc\u05d5\u05de\u05e8 \u05dc\u05d0 \u05e2\u05d5\u05de\u05e7 \u05e1\u05d5\u05e4\u05d9 \u05db\u05de\u05d5 \u05d1- TC0.\", Word.InsertLocation.replace);\nparas.items[7].insertText(\"\u05de\u05d3\u05d2\u05d9\u05dd \u05e9-SSMs \u05d0\u05d9\u05e0\u05dd \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05e2\u05e7\u05d5\u05d1 \u05d1\u05de\u05d3\u05d5\u05d9\u05e7 \u05d0\u05d7\u05e8 \u05de\u05d4\u05dc\u05db\u05d9 \u05e9\u05d7\u05de\u05d8, \u05dc\u05db\u05ea\u05d5\u05d1 \u05e7\u05d5\u05d3 \u05de\u05d5\u05e8\u05db\u05d1, \u05d0\u05d5 \u05dc\u05e2\u05e7\u05d5\u05d1 \u05d0\u05d7\u05e8 \u05d9\u05e9\u05d5\u05d9\u05d5\u05ea \u05d1\u05e0\u05e8\u05d8\u05d9\u05d1\u05d9\u05dd.\\u000b\", Word.InsertLocation.replace);\nparas.items[8].insertText(\"\u05d1\u05d3\u05d9\u05e7\u05d5\u05ea \u05d0\u05de\u05e4\u05d9\u05e8\u05d9\u05d5\u05ea \u05e9\u05d1\u05d5\u05e6\u05e2\u05d5 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05d4\u05de\u05d0\u05de\u05e8:\", Word.InsertLocation.replace);\nparas.items[9].insertText(\"\u05de\u05e1\u05e4\u05e7 \u05e8\u05d0\u05d9\u05d5\u05ea \u05e0\u05d9\u05e1\u05d9\u05d5\u05e0\u05d9\u05d5\u05ea \u05d4\u05de\u05e8\u05d0\u05d5\u05ea \u05e9-SSMs \u05d1\u05e1\u05d2\u05e0\u05d5\u05df Mamba \u05d5\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd \u05de\u05ea\u05e7\u05e9\u05d9\u05dd \u05d1\u05de\u05e9\u05d9\u05de\u05d5\u05ea permutation composition.\", Word.InsertLocation.replace);\nparas.items[10].insertText(\"\u05de\u05e8\u05d0\u05d4 \u05e9-SSMs \u05d3\u05d5\u05e8\u05e9\u05d9\u05dd \u05e2\u05d5\u05de\u05e7 \u05d2\u05d3\u05dc \u05db\u05d3\u05d9 \u05f4\u05dc\u05d8\u05e4\u05dc\u05f4 \u05d1\u05e8\u05e6\u05e4\u05d9\u05dd \u05d0\u05e8\u05d5\u05db\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05dc\u05de\u05d9\u05d3\u05d5\u05dc \u05e4\u05e2\u05d5\u05dc\u05d5\u05ea \u05e7\u05d1\u05d5\u05e6\u05d4 \u05f4\u05ea\u05de\u05d5\u05e8\u05ea\u05d9\u05d5\u05ea\u05f4\", Word.InsertLocation.replace);\nparas.items[11].insertText(\"\u05de\u05d3\u05d2\u05d9\u05dd \u05e9-RNNs \u05d1\u05e9\u05db\u05d1\u05d4 \u05d9\u05d7\u05d9\u05d3\u05d4 \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05e4\u05ea\u05d5\u05e8 \u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d0\u05dc\u05d5 \u05e9-SSMs \u05d0\u05d9\u05e0\u05dd \u05d9\u05db\u05d5\u05dc\u05d9\u05dd (\u05db\u05e0\u05e8\u05d0\u05d4 \u05d1\u05d2\u05dc\u05dc \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05d5\u05ea \u05d1\u05d9\u05df \u05d4\u05de\u05e2\u05d1\u05d9\u05e8\u05d9\u05dd \u05e9\u05dc \u05d4\u05de\u05e6\u05d1\u05d9\u05dd \u05d4\u05d7\u05d1\u05d5\u05d9\u05d9\u05dd \u05d1-SSMs).\\u000b\", Word.InsertLocation.replace);\nparas.items[12].insertText(\"\u05e9\u05db\u05dc\u05d5\u05dc\u05d9 \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e0\u05d9\u05d5\u05ea \u05d4\u05de\u05d5\u05e6\u05e2\u05d9\u05dd \u05d1\u05de\u05d0\u05de\u05e8:\", Word.InsertLocation.replace);\nparas.items[13].insertText(\"\u05de\u05e6\u05d9\u05e2 2 \u05d3\u05e8\u05db\u05d9\u05dd \u05dc\u05d4\u05e8\u05d7\u05d9\u05d1 SSMs \u05de\u05e2\u05d1\u05e8 \u05dc\u05de\u05d2\u05d1\u05dc\u05d5\u05ea TC0: \u05d4\u05d5\u05e1\u05e4\u05ea \u05d0\u05d9-\u05dc\u05d9\u05e0\u05d9\u05d0\u05e8\u05d9\u05d5\u05ea (RNN-SSM)  \u05d5\u05d4\u05e4\u05d9\u05db\u05ea \u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea \u05d4\u05de\u05e2\u05d1\u05e8 \u05dc\u05ea\u05dc\u05d5\u05d9\u05d5\u05ea \u05d1\u05e7\u05dc\u05d8 (WFA-SSM) - \u05e9\u05db\u05dc\u05d5\u05dc \u05e9\u05dc \u05de\u05de\u05d1\u05d4 \u05d4\u05de\u05d5\u05e1\u05d9\u05e3 \u05d0\u05d9 \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05d5\u05ea \u05dc\u05de\u05d8\u05e8\u05d9\u05e6\u05d4 A \u05e9\u05e0\u05d5\u05ea\u05e8\u05d4 \u05e7\u05d1\u05d5\u05e2\u05d4 \u05d1\u05de\u05de\u05d1\u05d4.\", Word.InsertLocation.replace);\nparas.items[14].insertText(\"\u05d4\u05e9\u05e4\u05e2\u05d4 \u05d5\u05d4\u05e9\u05dc\u05db\u05d5\u05ea \u05e9\u05dc \u05d4\u05de\u05d0\u05de\u05e8:\", Word.InsertLocation.replace);\nparas.items[15].insertText(\"\u05de\u05d0\u05ea\u05d2\u05e8 \u05d4\u05e0\u05d7\u05d5\u05ea \u05dc\u05d2\u05d1\u05d9 \u05d9\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea SSMs \u05e2\u05dc \u05e4\u05e0\u05d9 \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd\", Word.InsertLocation.replace);\nparas.items[16].insertText(\"\u05de\u05e6\u05d1\u05d9\u05e2 \u05e2\u05dc \u05d2\u05d9\u05e9\u05d5\u05ea \u05d4\u05d9\u05d1\u05e8\u05d9\u05d3\u05d9\u05d5\u05ea \u05e4\u05d5\u05d8\u05e0\u05e6\u05d9\u05d0\u05dc\u05d9\u05d5\u05ea \u05d4\u05de\u05e9\u05dc\u05d1\u05d5\u05ea \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d5\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea\", Word.InsertLocation.replace);\nparas.items[17].insertText(\"\u05e4\u05d5\u05ea\u05d7 \u05db\u05d9\u05d5\u05d5\u05e0\u05d9\u05dd \u05d7\u05d3\u05e9\u05d9\u05dd \u05dc\u05e4\u05d9\u05ea\u05d5\u05d7 \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d5\u05ea \u05e2\u05dd \u05d9\u05db\u05d5\u05dc\u05ea \u05d1\u05d9\u05d8\u05d5\u05d9 \u05de\u05e9\u05d5\u05e4\u05e8\u05ea \u05dc\u05d9\u05d9\u05e9\u05d5\u05de\u05d9 \u05e2\u05d9\u05d1\u05d5\u05d3 \u05e9\u05e4\u05d4 \u05d8\u05d1\u05e2\u05d9\u05ea \u05d5\u05e2\u05d1\u05d5\u05e8 \u05d3\u05d5\u05de\u05d9\u05d9\u05e0\u05d9\u05dd \u05e0\u05d5\u05e1\u05e4\u05d9\u05dd\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Insert the brand-new paragraph right after paragraph index 17\nparas.items[17].insertParagraph(\"\u05de\u05d3\u05d2\u05d9\u05e9 \u05d0\u05ea \u05d7\u05e9\u05d9\u05d1\u05d5\u05ea \u05d4\u05e0\u05d9\u05ea\u05d5\u05d7 \u05d4\u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9 \u05e9\u05dc \u05d4\u05ea\u05de\u05d0\u05ea \u05e9\u05dc \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05ea \u05de\u05d5\u05d3\u05dc \u05dc\u05de\u05e9\u05d9\u05de\u05d4 \u05e1\u05e4\u05e6\u05d9\u05e4\u05d9\u05ea \u05e9\u05d4\u05d5\u05d0 \u05de\u05ea\u05d5\u05db\u05e0\u05df \u05dc\u05e4\u05ea\u05d5\u05e8\", Word.InsertLocation.after);\nawait context.sync();\n\n// Paragraph \"\u05e1\u05d9\u05db\u05d5\u05dd:\" (index 18) is unchanged by the diff\n\n// The two trailing paragraphs (summary body + arxiv link), now at index 19/20\nparas.items[19].insertText(\"\u05de\u05d0\u05de\u05e8 \u05ea\u05d5\u05e8\u05dd \u05d4\u05df \u05de\u05d1\u05d7\u05d9\u05e0\u05d4 \u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9\u05ea \u05d5\u05d4\u05df \u05de\u05d1\u05d7\u05d9\u05e0\u05d4 \u05de\u05e2\u05e9\u05d9\u05ea \u05dc\u05d4\u05d1\u05e0\u05ea \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d5\u05ea \u05e9\u05dc \u05e8\u05e9\u05ea\u05d5\u05ea \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd. \u05d4\u05e0\u05d9\u05ea\u05d5\u05d7 \u05d4\u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9 \u05d4\u05e7\u05e4\u05d3\u05e0\u05d9, \u05d1\u05e9\u05d9\u05dc\u05d5\u05d1 \u05e2\u05dd \u05e8\u05d0\u05d9\u05d5\u05ea \u05d0\u05de\u05e4\u05d9\u05e8\u05d9\u05d5\u05ea \u05ea\u05d5\u05de\u05db\u05d5\u05ea, \u05de\u05e1\u05e4\u05e7 \u05ea\u05d5\u05d1\u05e0\u05d5\u05ea \u05d7\u05e9\u05d5\u05d1\u05d5\u05ea \u05dc\u05d2\u05d1\u05d9 \u05d4\u05de\u05d2\u05d1\u05dc\u05d5\u05ea \u05d4\u05d1\u05e1\u05d9\u05e1\u05d9\u05d5\u05ea \u05e9\u05dc SSMs.. \u05d1\u05e2\u05d5\u05d3 \u05e9\u05d7\u05dc\u05e7 \u05de\u05d4\u05ea\u05d5\u05e6\u05d0\u05d5\u05ea \u05d4\u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9\u05d5\u05ea \u05de\u05e1\u05ea\u05de\u05db\u05d5\u05ea \u05e2\u05dc \u05d4\u05e0\u05d7\u05d5\u05ea \u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9\u05d5\u05ea \u05e9\u05dc \u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea, \u05d4\u05d4\u05e9\u05dc\u05db\u05d5\u05ea \u05d4\u05de\u05e2\u05e9\u05d9\u05d5\u05ea \u05e0\u05ea\u05de\u05db\u05d5\u05ea \u05d4\u05d9\u05d8\u05d1 \u05d1\u05e8\u05d0\u05d9\u05d5\u05ea \u05d0\u05de\u05e4\u05d9\u05e8\u05d9\u05d5\u05ea.\", Word.InsertLocation.replace);\nparas.items[20].insertText(\"https://arxiv.org/abs/2404.08819\", Word.InsertLocation.replace);\nawait context.sync();\n\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Set-ParaText($paraIndex, $text) {\n    $r = $d.Paragraphs($paraIndex).Range\n    $r.MoveEnd(1, -1) | Out-Null\n    $r.Text = $text\n}\n\n# Paragraph 1: date + title line (single run, line break in the middle)\n$p1 = '\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 -27.11.24: \u26a1\ufe0f\ud83d\ude80' + ([char]0x0B) + 'The Illusion of State in State-Space Models'\nSet-ParaText 1 $p1\n\nSet-ParaText 2 '\u05de\u05d0\u05de\u05e8 \u05d7\u05e9\u05d5\u05d1 \u05d6\u05d4 \u05d1\u05d5\u05d7\u05df \u05d0\u05ea \u05d4\u05de\u05d2\u05d1\u05dc\u05d5\u05ea \u05d4\u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9\u05d5\u05ea \u05e9\u05dc State Space Models \u05d0\u05d5 (SSMs), \u05d0\u05e9\u05e8 \u05e6\u05de\u05d7\u05d5 \u05db\u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d4 \u05d7\u05dc\u05d5\u05e4\u05d9\u05ea \u05dc\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd \u05e2\u05d1\u05d5\u05e8 \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05d2\u05d3\u05d5\u05dc\u05d9\u05dd. \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05d3\u05d2\u05d9\u05de\u05d9\u05dd \u05e9\u05dc\u05de\u05e8\u05d5\u05ea \u05e2\u05d9\u05e6\u05d5\u05d1\u05dd \u05e9\u05e0\u05e8\u05d0\u05d4 Recurrent \u05d5\u05d1\u05e2\u05dc \u05de\u05e6\u05d1 (\u05db\u05dc\u05d5\u05de\u05e8 stateful), \u05dc\u05de\u05e2\u05e9\u05d4 SSMs (\u05db\u05de\u05d5 \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd) \u05de\u05d5\u05d2\u05d1\u05dc\u05d9\u05dd \u05d1\u05d0\u05d5\u05e4\u05df \u05d1\u05e1\u05d9\u05e1\u05d9 \u05d1\u05d9\u05db\u05d5\u05dc\u05ea\u05dd \u05dc\u05d1\u05d8\u05d0 \u05d7\u05d9\u05e9\u05d5\u05d1 \"\u05e8\u05e6\u05d9\u05e3\", \u05de\u05db\u05d9\u05d5\u05d5\u05df \u05e9\u05d0\u05d9\u05e0\u05dd \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05d7\u05e9\u05d1 \u05d3\u05d1\u05e8 \u05de\u05d7\u05d5\u05e5 \u05dc\u05de\u05d7\u05dc\u05e7\u05ea \u05d4\u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea TC0. \u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05de\u05de\u05d7\u05dc\u05e7\u05ea TC0 \u05de\u05d5\u05d2\u05d3\u05e8\u05d5\u05ea \u05db\u05db\u05d0\u05dc\u05d5 \u05e9\u05e0\u05d9\u05ea\u05df \u05dc\u05d9\u05d9\u05e6\u05d2\u05df \u05e2\u05dd \u05e9\u05e8\u05e9\u05e8\u05d0\u05d5\u05ea \u05d1\u05d5\u05dc\u05d9\u05d0\u05e0\u05d9\u05d5\u05ea \u05d1\u05e1\u05d9\u05e1\u05d9\u05d5\u05ea (\u05d5\u05d7\u05d9\u05e9\u05d5\u05d1\u05d9 \u05e1\u05e3 \u05d5- majority vote) \u05d1\u05e2\u05d5\u05de\u05e7 \u05e1\u05d5\u05e4\u05d9 (\u05dc\u05de\u05e9\u05dc \u05d7\u05d9\u05d1\u05d5\u05e8 \u05e9\u05dc \u05de\u05e1\u05e4\u05e8\u05d9\u05dd, \u05de\u05db\u05e4\u05dc\u05d4 \u05d0\u05d5 \u05de\u05d9\u05d5\u05df \u05e9\u05dc n \u05de\u05e1\u05e4\u05e8\u05d9\u05dd). \u05de\u05d3\u05d5\u05d1\u05e8 \u05d1\u05de\u05d7\u05dc\u05e7\u05d4 \u05d4\u05db\u05d9 \"\u05e4\u05e9\u05d5\u05d8\u05d4\" \u05d1\u05d4\u05d9\u05e8\u05e8\u05db\u05d9\u05d4 \u05e9\u05dc \u05ea\u05d5\u05e8\u05d4 \u05e1\u05d9\u05d1\u05d5\u05db\u05d9\u05d5\u05ea circuit (\u05db\u05dc\u05d5\u05de\u05e8 circuit complexity).'\nSet-ParaText 3 '\u05de\u05e9\u05de\u05e2\u05d5\u05ea \u05d4\u05d3\u05d1\u05e8 \u05d4\u05d9\u05d0 \u05e9-SSMs \u05d0\u05d9\u05e0\u05dd \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05e4\u05ea\u05d5\u05e8 \u05d1\u05e2\u05d9\u05d5\u05ea \u05de\u05e1\u05d5\u05d2 permutation composition \u05e9- RNNs \u05d1\u05e2\u05dc\u05d5\u05ea \u05e9\u05db\u05d1\u05d4 \u05d0\u05d7\u05ea \u05de\u05e1\u05d5\u05d2\u05dc\u05d5\u05ea \u05dc\u05e4\u05ea\u05d5\u05e8.'\nSet-ParaText 4 '\u05ea\u05e8\u05d5\u05de\u05d5\u05ea \u05de\u05e8\u05db\u05d6\u05d9\u05d5\u05ea \u05e9\u05dc \u05d4\u05de\u05d0\u05de\u05e8:'\nSet-ParaText 5 '\u05e0\u05d9\u05ea\u05d5\u05d7 \u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9:'\nSet-ParaText 6 '\u05de\u05d5\u05db\u05d9\u05d7 \u05e9\u05d2\u05dd SSMs \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05d9\u05dd \u05d5\u05d2\u05dd SSMs \u05d1\u05e1\u05d2\u05e0\u05d5\u05df Mamba \u05de\u05d5\u05d2\u05d1\u05dc\u05d9\u05dd \u05dc\u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea \u05d7\u05d9\u05e9\u05d5\u05d1\u05d9\u05ea TC0'\nSet-ParaText 7 '\u05de\u05e8\u05d0\u05d4 \u05e9-SSMs \u05d0\u05d9\u05e0\u05dd \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05e4\u05ea\u05d5\u05e8 \u05d1\u05e2\u05d9\u05d5\u05ea \u05e9\u05dc\u05de\u05d5\u05ea-NC1 (\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05e9\u05e0\u05d9\u05ea\u05df \u05dc\u05d9\u05d9\u05e6\u05d2 \u05d0\u05d5\u05ea\u05df \u05e2\u05dd \u05e4\u05e2\u05d5\u05dc\u05d5\u05ea \u05d1\u05d5\u05dc\u05d9\u05d0\u05e0\u05d9\u05d5\u05ea \u05d1\u05e2\u05d5\u05de\u05e7 \u05dc\u05d5\u05d2\u05e8\u05d9\u05ea\u05de\u05d9 \u05de\u05de\u05d9\u05de\u05d3 \u05d4\u05d1\u05e2\u05d9\u05d4 - \u05de\u05e1\u05e4\u05e8 \u05de\u05e9\u05ea\u05e0\u05d9\u05dd \u05d1\u05d2\u05d3\u05d5\u05dc) \u05db\u05de\u05d5 \u05d4\u05e8\u05db\u05d1\u05ea \u05ea\u05de\u05d5\u05e8\u05d5\u05ea. \u05db\u05dc\u05d5\u05de\u05e8 \u05dc\u05d0 \u05e2\u05d5\u05de\u05e7 \u05e1\u05d5\u05e4\u05d9 \u05db\u05de\u05d5 \u05d1- TC0.'\n$p8 = '\u05de\u05d3\u05d2\u05d9\u05dd \u05e9-SSMs \u05d0\u05d9\u05e0\u05dd \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05e2\u05e7\u05d5\u05d1 \u05d1\u05de\u05d3\u05d5\u05d9\u05e7 \u05d0\u05d7\u05e8 \u05de\u05d4\u05dc\u05db\u05d9 \u05e9\u05d7\u05de\u05d8, \u05dc\u05db\u05ea\u05d5\u05d1 \u05e7\u05d5\u05d3 \u05de\u05d5\u05e8\u05db\u05d1, \u05d0\u05d5 \u05dc\u05e2\u05e7\u05d5\u05d1 \u05d0\u05d7\u05e8 \u05d9\u05e9\u05d5\u05d9\u05d5\u05ea \u05d1\u05e0\u05e8\u05d8\u05d9\u05d1\u05d9\u05dd.' + ([char]0x0B)\nSet-ParaText 8 $p8\nSet-ParaText 9 '\u05d1\u05d3\u05d9\u05e7\u05d5\u05ea \u05d0\u05de\u05e4\u05d9\u05e8\u05d9\u05d5\u05ea \u05e9\u05d1\u05d5\u05e6\u05e2\u05d5 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05d4\u05de\u05d0\u05de\u05e8:'\nSet-ParaText 10 '\u05de\u05e1\u05e4\u05e7 \u05e8\u05d0\u05d9\u05d5\u05ea \u05e0\u05d9\u05e1\u05d9\u05d5\u05e0\u05d9\u05d5\u05ea \u05d4\u05de\u05e8\u05d0\u05d5\u05ea \u05e9-SSMs \u05d1\u05e1\u05d2\u05e0\u05d5\u05df Mamba \u05d5\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd \u05de\u05ea\u05e7\u05e9\u05d9\u05dd \u05d1\u05de\u05e9\u05d9\u05de\u05d5\u05ea permutation composition.'\nSet-ParaText 11 '\u05de\u05e8\u05d0\u05d4 \u05e9-SSMs \u05d3\u05d5\u05e8\u05e9\u05d9\u05dd \u05e2\u05d5\u05de\u05e7 \u05d2\u05d3\u05dc \u05db\u05d3\u05d9 \u05f4\u05dc\u05d8\u05e4\u05dc\u05f4 \u05d1\u05e8\u05e6\u05e4\u05d9\u05dd \u05d0\u05e8\u05d5\u05db\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05dc\u05de\u05d9\u05d3\u05d5\u05dc \u05e4\u05e2\u05d5\u05dc\u05d5\u05ea \u05e7\u05d1\u05d5\u05e6\u05d4 \u05f4\u05ea\u05de\u05d5\u05e8\u05ea\u05d9\u05d5\u05ea\u05f4'\n$p12 = '\u05de\u05d3\u05d2\u05d9\u05dd \u05e9-RNNs \u05d1\u05e9\u05db\u05d1\u05d4 \u05d9\u05d7\u05d9\u05d3\u05d4 \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05e4\u05ea\u05d5\u05e8 \u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d0\u05dc\u05d5 \u05e9-SSMs \u05d0\u05d9\u05e0\u05dd \u05d9\u05db\u05d5\u05dc\u05d9\u05dd (\u05db\u05e0\u05e8\u05d0\u05d4 \u05d1\u05d2\u05dc\u05dc \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05d5\u05ea \u05d1\u05d9\u05df \u05d4\u05de\u05e2\u05d1\u05d9\u05e8\u05d9\u05dd \u05e9\u05dc \u05d4\u05de\u05e6\u05d1\u05d9\u05dd \u05d4\u05d7\u05d1\u05d5\u05d9\u05d9\u05dd \u05d1-SSMs).' + ([char]0x0B)\nSet-ParaText 12 $p12\nSet-ParaText 13 '\u05e9\u05db\u05dc\u05d5\u05dc\u05d9 \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e0\u05d9\u05d5\u05ea \u05d4\u05de\u05d5\u05e6\u05e2\u05d9\u05dd \u05d1\u05de\u05d0\u05de\u05e8:'\nSet-ParaText 14 '\u05de\u05e6\u05d9\u05e2 2 \u05d3\u05e8\u05db\u05d9\u05dd \u05dc\u05d4\u05e8\u05d7\u05d9\u05d1 SSMs \u05de\u05e2\u05d1\u05e8 \u05dc\u05de\u05d2\u05d1\u05dc\u05d5\u05ea TC0: \u05d4\u05d5\u05e1\u05e4\u05ea \u05d0\u05d9-\u05dc\u05d9\u05e0\u05d9\u05d0\u05e8\u05d9\u05d5\u05ea (RNN-SSM)  \u05d5\u05d4\u05e4\u05d9\u05db\u05ea \u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea \u05d4\u05de\u05e2\u05d1\u05e8 \u05dc\u05ea\u05dc\u05d5\u05d9\u05d5\u05ea \u05d1\u05e7\u05dc\u05d8 (WFA-SSM) - \u05e9\u05db\u05dc\u05d5\u05dc \u05e9\u05dc \u05de\u05de\u05d1\u05d4 \u05d4\u05de\u05d5\u05e1\u05d9\u05e3 \u05d0\u05d9 \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05d5\u05ea \u05dc\u05de\u05d8\u05e8\u05d9\u05e6\u05d4 A \u05e9\u05e0\u05d5\u05ea\u05e8\u05d4 \u05e7\u05d1\u05d5\u05e2\u05d4 \u05d1\u05de\u05de\u05d1\u05d4.'\nSet-ParaText 15 '\u05d4\u05e9\u05e4\u05e2\u05d4 \u05d5\u05d4\u05e9\u05dc\u05db\u05d5\u05ea \u05e9\u05dc \u05d4\u05de\u05d0\u05de\u05e8:'\nSet-ParaText 16 '\u05de\u05d0\u05ea\u05d2\u05e8 \u05d4\u05e0\u05d7\u05d5\u05ea \u05dc\u05d2\u05d1\u05d9 \u05d9\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea SSMs \u05e2\u05dc \u05e4\u05e0\u05d9 \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd'\nSet-ParaText 17 '\u05de\u05e6\u05d1\u05d9\u05e2 \u05e2\u05dc \u05d2\u05d9\u05e9\u05d5\u05ea \u05d4\u05d9\u05d1\u05e8\u05d9\u05d3\u05d9\u05d5\u05ea \u05e4\u05d5\u05d8\u05e0\u05e6\u05d9\u05d0\u05dc\u05d9\u05d5\u05ea \u05d4\u05de\u05e9\u05dc\u05d1\u05d5\u05ea \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d5\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea'\nSet-ParaText 18 '\u05e4\u05d5\u05ea\u05d7 \u05db\u05d9\u05d5\u05d5\u05e0\u05d9\u05dd \u05d7\u05d3\u05e9\u05d9\u05dd \u05dc\u05e4\u05d9\u05ea\u05d5\u05d7 \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d5\u05ea \u05e2\u05dd \u05d9\u05db\u05d5\u05dc\u05ea \u05d1\u05d9\u05d8\u05d5\u05d9 \u05de\u05e9\u05d5\u05e4\u05e8\u05ea \u05dc\u05d9\u05d9\u05e9\u05d5\u05de\u05d9 \u05e2\u05d9\u05d1\u05d5\u05d3 \u05e9\u05e4\u05d4 \u05d8\u05d1\u05e2\u05d9\u05ea \u05d5\u05e2\u05d1\u05d5\u05e8 \u05d3\u05d5\u05de\u05d9\u05d9\u05e0\u05d9\u05dd \u05e0\u05d5\u05e1\u05e4\u05d9\u05dd'\n\n# Insert the brand-new paragraph after the (now-updated) paragraph 18\n$d.Paragraphs(18).Range.InsertParagraphAfter() | Out-Null\nSet-ParaText 19 '\u05de\u05d3\u05d2\u05d9\u05e9 \u05d0\u05ea \u05d7\u05e9\u05d9\u05d1\u05d5\u05ea \u05d4\u05e0\u05d9\u05ea\u05d5\u05d7 \u05d4\u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9 \u05e9\u05dc \u05d4\u05ea\u05de\u05d0\u05ea \u05e9\u05dc \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05ea \u05de\u05d5\u05d3\u05dc \u05dc\u05de\u05e9\u05d9\u05de\u05d4 \u05e1\u05e4\u05e6\u05d9\u05e4\u05d9\u05ea \u05e9\u05d4\u05d5\u05d0 \u05de\u05ea\u05d5\u05db\u05e0\u05df \u05dc\u05e4\u05ea\u05d5\u05e8'\n\n# Paragraph 20 (\"\u05e1\u05d9\u05db\u05d5\u05dd:\") is unchanged by the diff\n\n# Paragraph 21: summary body text\nSet-ParaText 21 '\u05de\u05d0\u05de\u05e8 \u05ea\u05d5\u05e8\u05dd \u05d4\u05df \u05de\u05d1\u05d7\u05d9\u05e0\u05d4 \u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9\u05ea \u05d5\u05d4\u05df \u05de\u05d1\u05d7\u05d9\u05e0\u05d4 \u05de\u05e2\u05e9\u05d9\u05ea \u05dc\u05d4\u05d1\u05e0\u05ea \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d5\u05ea \u05e9\u05dc \u05e8\u05e9\u05ea\u05d5\u05ea \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd. \u05d4\u05e0\u05d9\u05ea\u05d5\u05d7 \u05d4\u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9 \u05d4\u05e7\u05e4\u05d3\u05e0\u05d9, \u05d1\u05e9\u05d9\u05dc\u05d5\u05d1 \u05e2\u05dd \u05e8\u05d0\u05d9\u05d5\u05ea \u05d0\u05de\u05e4\u05d9\u05e8\u05d9\u05d5\u05ea \u05ea\u05d5\u05de\u05db\u05d5\u05ea, \u05de\u05e1\u05e4\u05e7 \u05ea\u05d5\u05d1\u05e0\u05d5\u05ea \u05d7\u05e9\u05d5\u05d1\u05d5\u05ea \u05dc\u05d2\u05d1\u05d9 \u05d4\u05de\u05d2\u05d1\u05dc\u05d5\u05ea \u05d4\u05d1\u05e1\u05d9\u05e1\u05d9\u05d5\u05ea \u05e9\u05dc SSMs.. \u05d1\u05e2\u05d5\u05d3 \u05e9\u05d7\u05dc\u05e7 \u05de\u05d4\u05ea\u05d5\u05e6\u05d0\u05d5\u05ea \u05d4\u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9\u05d5\u05ea \u05de\u05e1\u05ea\u05de\u05db\u05d5\u05ea \u05e2\u05dc \u05d4\u05e0\u05d7\u05d5\u05ea \u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9\u05d5\u05ea \u05e9\u05dc \u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea, \u05d4\u05d4\u05e9\u05dc\u05db\u05d5\u05ea \u05d4\u05de\u05e2\u05e9\u05d9\u05d5\u05ea \u05e0\u05ea\u05de\u05db\u05d5\u05ea \u05d4\u05d9\u05d8\u05d1 \u05d1\u05e8\u05d0\u05d9\u05d5\u05ea \u05d0\u05de\u05e4\u05d9\u05e8\u05d9\u05d5\u05ea.'\n\n# Paragraph 22: arxiv link\nSet-ParaText 22 'https://arxiv.org/abs/2404.08819'\n\n"}
